# DO July 2024 + EFO 3.68 release
# Update source_version for Disease Ontology (row 3) and Experimental Factor
# Ontology (row 4) on the phen_oncox sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Disease Ontology source_version: v2024-05-29 -> v2024-07-31
$ws.Range("E3").Value = "v2024-07-31"

# Experimental Factor Ontology source_version: v3.67.0 -> v3.68.0
$ws.Range("E4").Value = "v3.68.0"
